$d = $word.ActiveDocument

# --- Paragraph 1: text + paragraph formatting updates ---
$p1 = $d.Paragraphs(1)

# Replace the paragraph's visible text (everything up to, but excluding,
# the paragraph mark) with the new ID placeholder. This also collapses
# the old trailing-space run into the single resulting run.
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$r1.Text = "**ID__AFFARS_SUBPART_5328_1__ID**"

# Re-fetch the (now single-run) paragraph and update its indentation.
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Add a paragraph border (top/left/bottom/right) with 5pt space from text.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

Write-Output "done"
